$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TTLE")

# --- 1. Duplicate the last data row (row 7, "motorbikes") six times so rows
#        8-13 inherit the same cell styles (s="4" on B/C) as rows 2-7. ---
$ws.Rows(7).Copy()
$ws.Rows(8).Insert()
$ws.Rows(7).Copy()
$ws.Rows(9).Insert()
$ws.Rows(7).Copy()
$ws.Rows(10).Insert()
$ws.Rows(7).Copy()
$ws.Rows(11).Insert()
$ws.Rows(7).Copy()
$ws.Rows(12).Insert()
$ws.Rows(7).Copy()
$ws.Rows(13).Insert()

# --- 2. Relabel the passenger rows (2-7), prefixing the existing mode name. ---
$ws.Range("A2").Value = "passenger LDVs"
$ws.Range("A3").Value = "passenger HDVs"
$ws.Range("A4").Value = "passenger aircraft"
$ws.Range("A5").Value = "passenger rail"
$ws.Range("A6").Value = "passenger ships"
$ws.Range("A7").Value = "passenger motorbikes"

# --- 3. Label the new freight rows (8-13), same mode order as above. ---
$ws.Range("A8").Value = "freight LDVs"
$ws.Range("B8").Value = -3
$ws.Range("A9").Value = "freight HDVs"
$ws.Range("B9").Value = -3
$ws.Range("A10").Value = "freight aircraft"
$ws.Range("B10").Value = -3
$ws.Range("A11").Value = "freight rail"
$ws.Range("B11").Value = -3
$ws.Range("A12").Value = "freight ships"
$ws.Range("B12").Value = -3
$ws.Range("A13").Value = "freight motorbikes"
$ws.Range("B13").Value = -3

# --- 4. Header row: column B now holds a single "Value" column (was
#        "Passenger"); column C's "Freight" header goes away along with the
#        now-unneeded Freight data column. ---
$ws.Range("B1").Value = "Value"
$ws.Range("C1").ClearContents()
$ws.Range("C2:C13").Clear()

# --- 5. Widen column A to fit the longer "passenger/freight ..." labels. ---
$ws.Columns.Item(1).ColumnWidth = 21.7
